$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as text to match the source data format)
$updates = @{
    'D2' = '295.32'
    'E2' = '-2.04%'
    'D3' = '31.16'
    'E3' = '-2.62%'
    'D4' = '4.961'
    'E4' = '-1.07%'
    'D5' = '0.07366'
    'E5' = '-6.32%'
    'D6' = '1.843'
    'E6' = '-13.02%'
    'D7' = '7.680'
    'E7' = '-1.54%'
    'D8' = '3.751'
    'E8' = '-0.82%'
    'D9' = '0.9110'
    'E9' = '-1.70%'
    'D10' = '0.1672'
    'E10' = '-4.09%'
    'D11' = '0.07618'
    'E11' = '-3.87%'
    'D12' = '0.08127'
    'E12' = '-7.03%'
    'D13' = '0.02985'
    'E13' = '-4.39%'
    'D14' = '0.1001'
    'E14' = '0.11%'
    'D15' = '0.001496'
    'E15' = '-0.82%'
    'D16' = '0.005703'
    'E16' = '-5.04%'
    'E17' = '0.26%'
    'D18' = '2.098'
    'E18' = '-7.56%'
    'D19' = '0.3272'
    'E19' = '-0.06%'
    'D20' = '0.1305'
    'D21' = '4.337'
    'E21' = '4.61%'
    'D22' = '0.1998'
    'E22' = '11.54%'
    'D23' = '0.04473'
    'E23' = '-2.71%'
    'D24' = '0.001225'
    'E24' = '-0.77%'
    'D25' = '0.004057'
    'E25' = '-9.61%'
    'D26' = '0.0001250'
    'E26' = '0.17%'
    'D39' = '0.01665'
    'E39' = '-4.08%'
    'D40' = '0.04426'
    'E40' = '-6.85%'
    'D41' = '0.007409'
    'E41' = '1.97%'
    'D42' = '0.1326'
    'E42' = '-2.51%'
    'D43' = '0.002055'
    'E43' = '-1.05%'
    'D44' = '0.01109'
    'E44' = '3.57%'
    'D45' = '0.00005974'
    'E45' = '-1.25%'
    'E46' = '0.16%'
    'D47' = '2.102'
    'E47' = '156.28%'
    'E48' = '-11.44%'
    'D49' = '0.00002099'
    'E49' = '0.16%'
    'E50' = '0.16%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format so numeric-looking / percent-looking strings are not
    # auto-converted into numbers by Excel, preserving the original text representation.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
